$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.001916437105737102
$ws.Range("E2").Value = 0.001916437105737102

$ws.Range("D3").Value = [double]"1.297430491786452E-06"
$ws.Range("E3").Value = [double]"1.297430491786452E-06"

$ws.Range("D4").Value = [double]"1.213798302970531E-16"
$ws.Range("E4").Value = [double]"1.213798302970531E-16"

$ws.Range("D5").Value = [double]"1.292785139986398E-14"
$ws.Range("E5").Value = [double]"1.292785139986398E-14"

$ws.Range("D6").Value = 0.01067046432653596
$ws.Range("E6").Value = 0.01067046432653596

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0

$ws.Range("D8").Value = 0.9346328624543723
$ws.Range("E8").Value = 0.06536713754562773

$ws.Range("D9").Value = 0.9999990856797381
$ws.Range("E9").Value = [double]"9.143202619465995E-07"

$ws.Range("D11").Value = 0.9999778045970797
$ws.Range("E11").Value = [double]"2.219540292025535E-05"
$ws.Range("F11").Value = 0.008027281612157822
